$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "CreatedAt: 2026-02-04T19:07:13"

$ws.Range("W4").Value = 226.03
$ws.Range("X4").Value = 199.96
$ws.Range("Y4").Value = 352.09
$ws.Range("Z4").Value = 202.85
$ws.Range("W6").Value = -10.17
$ws.Range("X6").Value = -9.4
$ws.Range("Y6").Value = -12.68
$ws.Range("Z6").Value = -6.29
$ws.Range("W9").Value = 217.3
$ws.Range("X9").Value = 197.69
$ws.Range("Y9").Value = 353.11
$ws.Range("Z9").Value = 207.89
$ws.Range("W11").Value = -18.9
$ws.Range("X11").Value = -11.66
$ws.Range("Y11").Value = -11.65
$ws.Range("Z11").Value = -1.25
$ws.Range("W14").Value = 217.3
$ws.Range("X14").Value = 197.69
$ws.Range("Y14").Value = 353.11
$ws.Range("Z14").Value = 208.1
$ws.Range("W16").Value = -18.9
$ws.Range("X16").Value = -11.66
$ws.Range("Y16").Value = -11.65
$ws.Range("Z16").Value = -1.04
$ws.Range("W19").Value = 125.85
$ws.Range("X19").Value = 120.65
$ws.Range("Y19").Value = 109.52
$ws.Range("Z19").Value = 112.55
$ws.Range("W20").Value = -96.77
$ws.Range("X20").Value = -76.48
$ws.Range("Y20").Value = -238.21
$ws.Range("Z20").Value = -87.97
$ws.Range("W21").Value = -13.58
$ws.Range("X21").Value = -12.22
$ws.Range("Y21").Value = -17.04
$ws.Range("Z21").Value = -8.619999999999999
$ws.Range("W24").Value = 222.62
$ws.Range("X24").Value = 197.13
$ws.Range("Y24").Value = 347.73
$ws.Range("Z24").Value = 200.52
$ws.Range("W26").Value = -13.58
$ws.Range("X26").Value = -12.22
$ws.Range("Y26").Value = -17.04
$ws.Range("Z26").Value = -8.619999999999999
$ws.Range("X29").Value = 117.72
$ws.Range("Y29").Value = 103.66
$ws.Range("Z29").Value = 109.15
$ws.Range("W30").Value = -96.77
$ws.Range("X30").Value = -76.48
$ws.Range("Y30").Value = -238.21
$ws.Range("Z30").Value = -87.97
$ws.Range("W31").Value = -17.29
$ws.Range("X31").Value = -15.15
$ws.Range("Y31").Value = -22.9
$ws.Range("Z31").Value = -12.02
$ws.Range("X34").Value = 180.8
$ws.Range("Y34").Value = 359.38
$ws.Range("Z34").Value = 215.16
$ws.Range("W36").Value = -20.89
$ws.Range("X36").Value = -11.1
$ws.Range("Y36").Value = -5.39
$ws.Range("Z36").Value = 6.02
$ws.Range("W37").Value = -14.32
$ws.Range("X37").Value = -17.45
$ws.Range("W39").Value = 226.03
$ws.Range("X39").Value = 199.96
$ws.Range("Y39").Value = 352.09
$ws.Range("Z39").Value = 202.85
$ws.Range("W41").Value = -10.17
$ws.Range("X41").Value = -9.4
$ws.Range("Y41").Value = -12.68
$ws.Range("Z41").Value = -6.29
$ws.Range("W44").Value = 249.68
$ws.Range("X44").Value = 219.68
$ws.Range("Y44").Value = 380.76
$ws.Range("Z44").Value = 218.54
$ws.Range("W46").Value = 13.48
$ws.Range("X46").Value = 10.33
$ws.Range("Y46").Value = 15.99
$ws.Range("Z46").Value = 9.4
$ws.Range("W49").Value = 228.44
$ws.Range("X49").Value = 204.25
$ws.Range("Y49").Value = 379.17
$ws.Range("Z49").Value = 215.61
$ws.Range("W51").Value = -7.77
$ws.Range("X51").Value = -5.11
$ws.Range("Y51").Value = 14.41
$ws.Range("Z51").Value = 6.47
$ws.Range("W54").Value = 226.9
$ws.Range("X54").Value = 209.15
$ws.Range("Y54").Value = 366.23
$ws.Range("Z54").Value = 211.47
$ws.Range("W56").Value = -9.300000000000001
$ws.Range("X56").Value = -0.21
$ws.Range("Y56").Value = 1.46
$ws.Range("Z56").Value = 2.33
$ws.Range("W59").Value = 244.01
$ws.Range("X59").Value = 214.29
$ws.Range("Y59").Value = 374.12
$ws.Range("Z59").Value = 214.28
$ws.Range("W61").Value = 7.81
$ws.Range("X61").Value = 4.93
$ws.Range("Y61").Value = 9.35
$ws.Range("Z61").Value = 5.14
$ws.Range("W64").Value = 246.56
$ws.Range("X64").Value = 216.95
$ws.Range("Y64").Value = 378
$ws.Range("Z64").Value = 216.28
$ws.Range("W66").Value = 10.36
$ws.Range("X66").Value = 7.59
$ws.Range("Y66").Value = 13.23
$ws.Range("Z66").Value = 7.14
$ws.Range("W69").Value = 244.26
$ws.Range("X69").Value = 213.85
$ws.Range("Y69").Value = 374.5
$ws.Range("Z69").Value = 214.28
$ws.Range("W71").Value = 8.06
$ws.Range("X71").Value = 4.49
$ws.Range("Y71").Value = 9.74
$ws.Range("Z71").Value = 5.14
$ws.Range("W74").Value = 244.01
$ws.Range("X74").Value = 214.07
$ws.Range("Y74").Value = 374.5
$ws.Range("Z74").Value = 214.72
$ws.Range("W76").Value = 7.81
$ws.Range("X76").Value = 4.71
$ws.Range("Y76").Value = 9.74
$ws.Range("Z76").Value = 5.58
$ws.Range("W79").Value = 236.2
$ws.Range("X79").Value = 209.36
$ws.Range("Y79").Value = 364.77
$ws.Range("Z79").Value = 209.14
$ws.Range("W84").Value = 223.25
$ws.Range("X84").Value = 211.26
$ws.Range("Y84").Value = 370.32
$ws.Range("Z84").Value = 211.47
$ws.Range("W86").Value = -12.95
$ws.Range("X86").Value = 1.9
$ws.Range("Y86").Value = 5.55
$ws.Range("Z86").Value = 2.33
$ws.Range("W89").Value = 218.91
$ws.Range("X89").Value = 194.03
$ws.Range("Y89").Value = 341.86
$ws.Range("Z89").Value = 197.12
$ws.Range("W91").Value = -17.29
$ws.Range("X91").Value = -15.33
$ws.Range("Y91").Value = -22.9
$ws.Range("Z91").Value = -12.02
